# The deck's applied design ("Integral") is swapped back to the plain,
# default "Office Theme" colour scheme. PowerPoint keeps the 12 theme
# slots (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) on the
# presentation's ThemeColorScheme - reassigning them is exactly what the
# Design/Variants gallery does when a different theme is applied, so we
# push the stock "Office Theme" RGB values into every slot.

function ToBGR($r, $g, $b) {
    # PowerPoint's ColorFormat.RGB (like VBA's RGB() macro) stores the
    # value as a Long in 0x00BBGGRR order, so pack the bytes that way.
    return $b * 65536 + $g * 256 + $r
}

$p = $ppt.ActivePresentation
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

$tcs.Item(1).RGB  = ToBGR 0x00 0x00 0x00   # dk1
$tcs.Item(2).RGB  = ToBGR 0xFF 0xFF 0xFF   # lt1
$tcs.Item(3).RGB  = ToBGR 0x44 0x54 0x6A   # dk2
$tcs.Item(4).RGB  = ToBGR 0xE7 0xE6 0xE6   # lt2
$tcs.Item(5).RGB  = ToBGR 0x5B 0x9B 0xD5   # accent1
$tcs.Item(6).RGB  = ToBGR 0xED 0x7D 0x31   # accent2
$tcs.Item(7).RGB  = ToBGR 0xA5 0xA5 0xA5   # accent3
$tcs.Item(8).RGB  = ToBGR 0xFF 0xC0 0x00   # accent4
$tcs.Item(9).RGB  = ToBGR 0x44 0x72 0xC4   # accent5
$tcs.Item(10).RGB = ToBGR 0x70 0xAD 0x47   # accent6
$tcs.Item(11).RGB = ToBGR 0x05 0x63 0xC1   # hlink
$tcs.Item(12).RGB = ToBGR 0x95 0x4F 0x72   # folHlink
